$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$s.Shapes.Item("CaixaDeTexto 6").Delete()
